# Updates the crypto price/volume table with freshly scraped values.
# Only the "Price" (column D) and "Volume(1h)" (column E) cells that
# actually changed are touched; everything else is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "25.962.42"; E = "  -0.24%  " }
    @{ Row = 3; D = "1.627.49"; E = "  -1.01%  " }
    @{ Row = 4; D = $null; E = "  -0.14%  " }
    @{ Row = 5; D = "214.26"; E = "  -0.80%  " }
    @{ Row = 6; D = $null; E = "  -0.70%  " }
    @{ Row = 7; D = $null; E = "  -0.12%  " }
    @{ Row = 8; D = $null; E = "  -2.06%  " }
    @{ Row = 9; D = $null; E = "  -3.26%  " }
    @{ Row = 10; D = "18.48"; E = "  -5.55%  " }
    @{ Row = 11; D = $null; E = "  -1.09%  " }
    @{ Row = 12; D = "1.853.44"; E = "  -1.01%  " }
    @{ Row = 13; D = "1.628.56"; E = "  -1.19%  " }
    @{ Row = 14; D = $null; E = "  -2.16%  " }
    @{ Row = 15; D = $null; E = "  -3.11%  " }
    @{ Row = 16; D = "25.967.90"; E = "  -0.36%  " }
    @{ Row = 17; D = "0.0₃0740"; E = "  -3.10%  " }
    @{ Row = 18; D = "61.37"; E = "  -3.27%  " }
    @{ Row = 19; D = $null; E = "  -0.09%  " }
    @{ Row = 20; D = "192.92"; E = "  -0.72%  " }
    @{ Row = 21; D = $null; E = "  -2.40%  " }
    @{ Row = 22; D = $null; E = "  -3.52%  " }
    @{ Row = 23; D = "6.07"; E = "  -2.08%  " }
    @{ Row = 24; D = "0.134"; E = "  +1.17%  " }
    @{ Row = 25; D = "144.01"; E = "  +0.53%  " }
    @{ Row = 26; D = $null; E = "  +0.05%  " }
    @{ Row = 27; D = $null; E = "  -3.86%  " }
    @{ Row = 28; D = "6.74"; E = "  -2.07%  " }
    @{ Row = 29; D = "15.26"; E = "  -1.65%  " }
    @{ Row = 30; D = $null; E = "  -1.26%  " }
    @{ Row = 31; D = $null; E = "  -2.23%  " }
    @{ Row = 32; D = $null; E = "  -4.07%  " }
    @{ Row = 33; D = $null; E = "  -5.24%  " }
    @{ Row = 34; D = $null; E = "  -2.53%  " }
    @{ Row = 36; D = "1.125.93"; E = "  -0.45%  " }
    @{ Row = 37; D = $null; E = "  -5.78%  " }
    @{ Row = 38; D = $null; E = "  -1.56%  " }
    @{ Row = 39; D = $null; E = "  -3.29%  " }
    @{ Row = 40; D = $null; E = "  -2.23%  " }
    @{ Row = 41; D = "98.15"; E = "  -1.04%  " }
    @{ Row = 42; D = "1.764.24"; E = "  -0.95%  " }
    @{ Row = 43; D = $null; E = "  -4.35%  " }
    @{ Row = 44; D = $null; E = "  -5.28%  " }
    @{ Row = 45; D = $null; E = "  +1.99%  " }
    @{ Row = 46; D = "54.47"; E = "  -3.53%  " }
    @{ Row = 47; D = $null; E = "  -0.34%  " }
    @{ Row = 48; D = $null; E = "  -0.53%  " }
    @{ Row = 49; D = $null; E = "  +0.08%  " }
    @{ Row = 50; D = "7.47"; E = "  -3.90%  " }
    @{ Row = 51; D = "0.0927"; E = "  -2.51%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel keeps these numeric-looking
        # strings (e.g. "214.26", "18.48") as plain text instead of
        # silently converting them to floating point numbers, matching
        # the original inline-string cell type.
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
